$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# bill_detail table (column F): drop the "product_id" field.
# F5 was "product_id" -> becomes "quantity" (shifted up from F6)
$ws.Range("F5").Value = "quantity"
# F6 used to hold "quantity"; it is now empty (cell removed)
$ws.Range("F6").ClearContents()

# product table (column E): reference the id fields instead of the name fields
$ws.Range("E12").Value = "manufacturer_id"
$ws.Range("E16").Value = "type_id"

# update the last selected cell to reflect where editing ended
[void]$ws.Range("F21").Select()
